$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are pre-formatted as Text so the literal string is preserved (matches the
# source data, which stores these as plain/inline strings, not numbers).
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D14",
    "D18",
    "D20",
    "D21",
    "D22",
    "D26",
    "D27",
    "D28",
    "D32",
    "D33",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D46",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '35.432.32'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '1.886.30'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  -0.75%  '
$ws.Range("D5").Value = '246.14'
$ws.Range("E5").Value = '  -3.76%  '
$ws.Range("D6").Value = '0.691'
$ws.Range("E6").Value = '  -4.69%  '
$ws.Range("E7").Value = '  -0.82%  '
$ws.Range("D8").Value = '43.08'
$ws.Range("E8").Value = '  +2.13%  '
$ws.Range("D9").Value = '0.352'
$ws.Range("E9").Value = '  -3.97%  '
$ws.Range("D10").Value = '0.0737'
$ws.Range("E10").Value = '  -4.21%  '
$ws.Range("D11").Value = '0.0970'
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").Value = '2.158.95'
$ws.Range("D14").Value = '0.741'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("D16").Value = '1.875.88'
$ws.Range("E16").Value = '  -2.14%  '
$ws.Range("D17").Value = '35.396.66'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").Value = '73.55'
$ws.Range("E18").Value = '  -2.13%  '
$ws.Range("D19").Value = '0.0₃0822'
$ws.Range("E19").Value = '  -3.22%  '
$ws.Range("D20").Value = '245.05'
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").Value = '12.82'
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("D22").Value = '4.93'
$ws.Range("E22").Value = '  -4.63%  '
$ws.Range("E23").Value = '  -0.84%  '
$ws.Range("E24").Value = '  +3.18%  '
$ws.Range("E25").Value = '  -9.15%  '
$ws.Range("D26").Value = '165.89'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").Value = '8.48'
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").Value = '18.36'
$ws.Range("E28").Value = '  -2.43%  '
$ws.Range("E29").Value = '  -4.07%  '
$ws.Range("D30").Value = '4.128.47'
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  +2.99%  '
$ws.Range("D32").Value = '4.24'
$ws.Range("E32").Value = '  -2.96%  '
$ws.Range("D33").Value = '0.0578'
$ws.Range("E33").Value = '  -2.54%  '
$ws.Range("E34").Value = '  -2.16%  '
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("D36").Value = '1.76'
$ws.Range("E36").Value = '  -12.08%  '
$ws.Range("D37").Value = '0.853'
$ws.Range("E37").Value = '  -7.41%  '
$ws.Range("D38").Value = '1.97'
$ws.Range("E38").Value = '  -4.07%  '
$ws.Range("E39").Value = '  +4.96%  '
$ws.Range("D40").Value = '97.60'
$ws.Range("E40").Value = '  -2.89%  '
$ws.Range("D41").Value = '0.0217'
$ws.Range("E41").Value = '  -1.88%  '
$ws.Range("D42").Value = '17.03'
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("E43").Value = '  -4.53%  '
$ws.Range("D44").Value = '1.293.69'
$ws.Range("E44").Value = '  -3.77%  '
$ws.Range("E45").Value = '  -5.54%  '
$ws.Range("D46").Value = '0.0805'
$ws.Range("E46").Value = '  +5.97%  '
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("D48").Value = '2.73'
$ws.Range("E48").Value = '  -0.97%  '
$ws.Range("D49").Value = '12.00'
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("D50").Value = '43.19'
$ws.Range("E50").Value = '  -4.27%  '
$ws.Range("D51").Value = '6.24'
$ws.Range("E51").Value = '  -7.57%  '
